$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the depth value (B11): 4.5 -> 7.5 (6in focal length prep change)
$ws.Range("B11").Value = 7.5

# Update the active cell selection to B12
$ws.Range("B12").Select()
